# Apply cryptos list update (prices/volumes refreshed; a few coin rows swapped order)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.889.92"
$ws.Range("E2").Value = "  -3.55%  "

$ws.Range("D3").Value = "3.060.62"
$ws.Range("E3").Value = "  -2.42%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "'541.79"
$ws.Range("E5").Value = "  -3.76%  "

$ws.Range("D6").Value = "'136.09"
$ws.Range("E6").Value = "  -7.84%  "

$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").Value = "3.050.94"
$ws.Range("E8").Value = "  -2.15%  "

$ws.Range("D9").Value = "'0.488"
$ws.Range("E9").Value = "  -2.08%  "

$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'6.48"
$ws.Range("E10").Value = "  -6.55%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.156"
$ws.Range("E11").Value = "  -1.07%  "

$ws.Range("D12").Value = "'0.456"
$ws.Range("E12").Value = "  -1.26%  "

$ws.Range("D13").Value = "'34.55"
$ws.Range("E13").Value = "  -4.66%  "

$ws.Range("D14").Value = "'0.0000215"
$ws.Range("E14").Value = "  -2.73%  "

$ws.Range("D15").Value = "3.565.88"
$ws.Range("E15").Value = "  -1.92%  "

$ws.Range("D16").Value = "63.045.71"
$ws.Range("E16").Value = "  -3.29%  "

$ws.Range("E17").Value = "  -1.36%  "

$ws.Range("D18").Value = "3.071.57"
$ws.Range("E18").Value = "  -1.95%  "

$ws.Range("D19").Value = "'494.59"
$ws.Range("E19").Value = "  -4.68%  "

$ws.Range("D20").Value = "'6.58"
$ws.Range("E20").Value = "  -2.07%  "

$ws.Range("D21").Value = "'13.31"
$ws.Range("E21").Value = "  -3.81%  "

$ws.Range("D22").Value = "'0.694"
$ws.Range("E22").Value = "  -0.93%  "

$ws.Range("D23").Value = "'7.09"
$ws.Range("E23").Value = "  -4.25%  "

$ws.Range("D24").Value = "'77.15"
$ws.Range("E24").Value = "  -2.10%  "

$ws.Range("D25").Value = "'12.16"
$ws.Range("E25").Value = "  -4.47%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "'2.70"
$ws.Range("E27").Value = "  -3.02%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'8.24"
$ws.Range("E28").Value = "  -5.54%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("D30").Value = "'1.92"
$ws.Range("E30").Value = "  -10.19%  "

$ws.Range("D31").Value = "'26.12"
$ws.Range("E31").Value = "  -0.31%  "

$ws.Range("D32").Value = "'1.13"
$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("D33").Value = "'2.49"
$ws.Range("E33").Value = "  -7.23%  "

$ws.Range("D34").Value = "'59.40"
$ws.Range("E34").Value = "  +12.02%  "

$ws.Range("D35").Value = "'519.12"
$ws.Range("E35").Value = "  -8.40%  "

$ws.Range("D36").Value = "'5.88"
$ws.Range("E36").Value = "  -3.15%  "

$ws.Range("D37").Value = "'5.12"
$ws.Range("E37").Value = "  -6.88%  "

$ws.Range("D38").Value = "'0.0397"
$ws.Range("E38").Value = "  -8.64%  "

$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "3.038.51"
$ws.Range("E39").Value = "  -1.60%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.0785"
$ws.Range("E40").Value = "  -4.24%  "

$ws.Range("E41").Value = "  -3.38%  "

$ws.Range("D42").Value = "'8.05"
$ws.Range("E42").Value = "  -1.97%  "

$ws.Range("D43").Value = "'2.62"
$ws.Range("E43").Value = "  -8.74%  "

$ws.Range("D44").Value = "'0.253"
$ws.Range("E44").Value = "  -1.43%  "

$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "'121.69"
$ws.Range("E46").Value = "  +2.83%  "

$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.01"
$ws.Range("E47").Value = "  -8.16%  "

$ws.Range("D48").Value = "'24.09"
$ws.Range("E48").Value = "  -3.81%  "

$ws.Range("B49").Value = "CoreDAO"
$ws.Range("C49").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D49").Value = "'2.45"
$ws.Range("E49").Value = "  +56.46%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.106"
$ws.Range("E50").Value = "  -1.80%  "

$ws.Range("E51").Value = "  -5.34%  "
